# Generate Report for Handback
# Update the "Overview" sheet and per-language sheets with the latest
# handoff/handback timestamps for the file that was just handed back
# (6e69aeeb-2812-4203-93f3-f3326070ca5b.md).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-13 13:01:04"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-13 13:00:54"
$zhcn.Range("K2").Value = "2016-08-13 13:01:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-13 13:01:35"
